# ============================================================
# Add new "Knarot" section (heading, body paragraphs, references)
# after the "BILAGA 1 - Fridlysta arter" paragraph, and update the
# header date from 2023-09-13 to 2023-09-15.
# ============================================================

$d = $word.ActiveDocument

# ---- Phase 1: insert each new paragraph as plain text, one after
#      another, recording the absolute [start,end) character range
#      each paragraph's text occupies in the document. Formatting
#      (italics) is applied afterwards in Phase 2, using those
#      recorded absolute positions -- this sidesteps COM "current
#      typing format" state leaking across runs/paragraphs.
$anchor = $d.Paragraphs.Last.Range
$anchor.Collapse(0)

# --- paragraph 1 (Heading1) ---
$anchor.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Style = "Heading1"
$r = $p.Range
$r.Collapse(0)
$t0 = "Knärot – ekologi samt krav på livsmiljön"
$r.InsertAfter($t0)
$pEnd0 = $r.End
$pStart0 = $pEnd0 - $t0.Length
$r.Collapse(0)
$anchor = $r

# --- paragraph 2 (body) ---
$anchor.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Style = "Normal"
$r = $p.Range
$r.Collapse(0)
$t1 = "Knärot är fridlyst enligt 8 och 15 §§ artskyddsförordningen och klassad som sårbar (VU) enligt rödlistan 2020. Knärot är beroende av hög och jämn luftfuktighet i gamla, ostörda skogsmiljöer och är känslig för snabba förändringar av ljus-/vindförhållanden eller uttorkning. På grund av ett alltför intensivt skogsbruk har den minskat med 40 (25-50) % under de senaste 60 åren och i framtiden bedöms minskningstakten uppgå till 30 (20-40) %. Till följd av att arten har en dokumenterat högre minskningstakt iförhållande till sin generationstid än vad som tidigare varit känt (data från Riksskogstaxeringen) höjdes den till hotkategori sårbar (VU) i rödlistan 2020 (Artdatabanken, 2021)."
$r.InsertAfter($t1)
$pEnd1 = $r.End
$pStart1 = $pEnd1 - $t1.Length
$r.Collapse(0)
$anchor = $r

# --- paragraph 3 (body) ---
$anchor.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Style = "Normal"
$r = $p.Range
$r.Collapse(0)
$t2 = "Samuel Johnsons doktorsavhandling “Retention Forestry as a Conservation Measure for Boreal Forest Ground Vegetation“ (SLU, Uppsala 2014) visar att det krävs väl tilltagna skyddszoner för att knärotens växtplatser inte ska ta skada av skogsbruksåtgärder i intilliggande områden: “Study III shows that retention patches smaller than 0.5 ha do not lifeboat the sensitive forest herb G. repens, a species that depend on stable microclimatic conditions typical for intact forest stands.” Vidare “More sensitive forest species are not lifeboated in retention patches ranging from 0.05 to 0.5 ha (Papers II & III).”"
$r.InsertAfter($t2)
$pEnd2 = $r.End
$pStart2 = $pEnd2 - $t2.Length
$r.Collapse(0)
$anchor = $r

# --- paragraph 4 (body) ---
$anchor.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Style = "Normal"
$r = $p.Range
$r.Collapse(0)
$t3 = "Johnsons (2014) rekommendation på minst 50 meters breda skyddszoner runt knärotens växtplatser motsvarar en areal på 0,78 hektar, vilket ligger i linje med andra studier som gjorts på känsliga skogsarter: “In study III I also show that translocated specimens of G. repens survives well in mature forests at least 50 m from the nearest edge to an open area. Moreover, measures of temperature and humidity show that such distances from an open area is far enough to offer a microclimate that is more stable compared to what present in retention patches of around 0.1 ha. This means that the very centre of a circular patch with radius 50 m (equals a size of 0.78 ha) should offer conditions similar to interior forest and would perhaps be a suitable habitat for G. repens and similar species. Previous studies from both North America and Sweden have also concluded that patches between 0.5 and one ha are sufficient for preserving interior forest vegetation as well as sensitive lichens and bryophytes (de Graaf & Roberts 2009; Halpern et al. 2012; Rudolphi et al. 2014).”"
$r.InsertAfter($t3)
$pEnd3 = $r.End
$pStart3 = $pEnd3 - $t3.Length
$r.Collapse(0)
$anchor = $r

# --- paragraph 5 (body) ---
$anchor.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Style = "Normal"
$r = $p.Range
$r.Collapse(0)
$t4 = "En nyligen publicerad vetenskaplig uppsats av Koelmeijer m.fl. (2022) inkluderar orkidén knärots skyddsbehov. I uppsatsen berörs problemet med uttorkning för växter, bl.a. för knärot, ett problem som blivit accentuerat på grund av den pågående klimatförändringen och torra somrar, t.ex. den exceptionellt torra sommaren 2018. I uppsatsen undersöks områden med tre olika avstånd från kalhyggeskant med avseende på skydd bl.a. för knärot. Det första området har avstånd upp till 20 m från hyggeskant (Strong edge effect), det andra 20 – 40 m från hyggeskant (Weak edge effect) och det tredje avser större avstånd från hyggeskant, där kanteffekten anses vara försumbar (Interior). Ett resultat var att man fann stor eller mycket stor uttorkningseffekt på känsliga och rödlistade skogsarter vid de kortare avstånden till hyggeskant, medan effekt av uttorkning inte konstaterades på större avstånd (Interior). För orkidén knärot fann man en rik förekomst (upp till 0,06 dm2/m2) på stort avstånd från hyggeskant (Interior), medan förekomsten var liten eller närmast försumbar i de områden som klassificerades som Weak edge effect respektive Strong edge effect. Arbetet påpekar att de allt oftare förekommande torra somrarna ger ytterligare skäl att utöka skyddsavståndet från hyggen till den fuktkrävande arten knärot (Koelmeijer m.fl., 2022)."
$r.InsertAfter($t4)
$pEnd4 = $r.End
$pStart4 = $pEnd4 - $t4.Length
$r.Collapse(0)
$anchor = $r

# --- paragraph 6 (body) ---
$anchor.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Style = "Normal"
$r = $p.Range
$r.Collapse(0)
$t5 = "Även Skogsstyrelsens egen vägledning för hänsyn till knärot ligger i linje med ovanstående forskningsstudier. Av vägledningen framgår det att för med hög sannolikhet kunna bevara befintliga förekomster krävs relativt stora avsättningar av uppvuxen skog med slutet och relativt tätt kronskikt. Som riktlinje kan krävas ett avstånd på 50 meter in från brynet för att vidmakthålla ett fungerande mikroklimat. Detta innebär att fristående hänsynsytor för många arter (kärlväxter, lavar och mossor) kan behöva ha en area överstigande 0,8 hektar (cirkelyta med radien 50 meter = 0,78 hektar) för att bibehålla lokalklimatet. Även ganska små förändringar i form av förändrade ljus- och fuktighetsförhållanden, till exempel till följd av gallring, kan leda till att arten försvinner till följd av konkurrens med mera ljuskrävande och snabbväxande arter (Skogsstyrelsen, 2022)."
$r.InsertAfter($t5)
$pEnd5 = $r.End
$pStart5 = $pEnd5 - $t5.Length
$r.Collapse(0)
$anchor = $r

# --- paragraph 7 (Heading2) ---
$anchor.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Style = "Heading2"
$r = $p.Range
$r.Collapse(0)
$t6 = "Referenser - knärot"
$r.InsertAfter($t6)
$pEnd6 = $r.End
$pStart6 = $pEnd6 - $t6.Length
$r.Collapse(0)
$anchor = $r

# --- paragraph 8 (body) ---
$anchor.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Style = "Normal"
$r = $p.Range
$r.Collapse(0)
$t7 = "de Graaf M & Roberts M.R., 2009. Short-term response of the herbaceous layer within leave patches after harvest. Forest Ecology and Management 257, 1014-1025"
$r.InsertAfter($t7)
$pEnd7 = $r.End
$pStart7 = $pEnd7 - $t7.Length
$r.Collapse(0)
$anchor = $r

# --- paragraph 9 (body) ---
$anchor.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Style = "Normal"
$r = $p.Range
$r.Collapse(0)
$t8 = "Halpern, C. B., Halaj, J., Evans, S. A., & Dovciak, M., 2012. Level and pattern of overstory retention interact to shape long-term responses of understories to timber harvest. Ecological Applications, 22, 2049-2064 "
$r.InsertAfter($t8)
$pEnd8 = $r.End
$pStart8 = $pEnd8 - $t8.Length
$r.Collapse(0)
$anchor = $r

# --- paragraph 10 (body) ---
$anchor.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Style = "Normal"
$r = $p.Range
$r.Collapse(0)
$t9 = "Koelmeijer, I. A., Ehrlén, J., Jönsson, M., De Frenne, P., Berg, P., Andersson, J., Weibull, H. & Hylander, N. 2022. Interactive effects of drought and edge exposure on old-growth forest understory species. Landscape Ecology, 37, sid 1839-1853"
$r.InsertAfter($t9)
$pEnd9 = $r.End
$pStart9 = $pEnd9 - $t9.Length
$r.Collapse(0)
$anchor = $r

# --- paragraph 11 (body) ---
$anchor.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Style = "Normal"
$r = $p.Range
$r.Collapse(0)
$t10 = "Rudolphi, J., Jönsson, M. T., & Gustafsson, L., 2014. Biological legacies buffer local species extinction after logging. Journal of Applied Ecology. 51, 53-62."
$r.InsertAfter($t10)
$pEnd10 = $r.End
$pStart10 = $pEnd10 - $t10.Length
$r.Collapse(0)
$anchor = $r

# --- paragraph 12 (body) ---
$anchor.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Style = "Normal"
$r = $p.Range
$r.Collapse(0)
$t11 = "Skogsstyrelsen, 2022. Vägledning för hänsyn till knärot. https://www.skogsstyrelsen.se/lag-och-tillsyn/artskydd/vagledningar-och-kunskapsstod-artskydd/vagledning-for-hansyn-till-knarot/"
$r.InsertAfter($t11)
$pEnd11 = $r.End
$pStart11 = $pEnd11 - $t11.Length
$r.Collapse(0)
$anchor = $r

# --- paragraph 13 (body) ---
$anchor.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Style = "Normal"
$r = $p.Range
$r.Collapse(0)
$t12 = "SLU Artdatabanken, 2021. Artfaktablad. Naturvård – artfakta. SLU Artdatabanken, Uppsala "
$r.InsertAfter($t12)
$pEnd12 = $r.End
$pStart12 = $pEnd12 - $t12.Length
$r.Collapse(0)
$anchor = $r

# ---- Phase 2: layer italic formatting onto the sub-ranges that need it.

# paragraph 3
$sStart = $pStart2 + 34
$sEnd = $pStart2 + 116
$sub = $d.Range($sStart, $sEnd)
$sub.Font.Italic = $true
$sStart = $pStart2 + 278
$sEnd = $pStart2 + 483
$sub = $d.Range($sStart, $sEnd)
$sub.Font.Italic = $true
$sStart = $pStart2 + 490
$sEnd = $pStart2 + 608
$sub = $d.Range($sStart, $sEnd)
$sub.Font.Italic = $true

# paragraph 4
$sStart = $pStart3 + 205
$sEnd = $pStart3 + 1070
$sub = $d.Range($sStart, $sEnd)
$sub.Font.Italic = $true

# paragraph 8
$sStart = $pStart7 + 33
$sEnd = $pStart7 + 113
$sub = $d.Range($sStart, $sEnd)
$sub.Font.Italic = $true

# paragraph 9
$sStart = $pStart8 + 62
$sEnd = $pStart8 + 176
$sub = $d.Range($sStart, $sEnd)
$sub.Font.Italic = $true

# paragraph 10
$sStart = $pStart9 + 117
$sEnd = $pStart9 + 207
$sub = $d.Range($sStart, $sEnd)
$sub.Font.Italic = $true

# paragraph 11
$sStart = $pStart10 + 54
$sEnd = $pStart10 + 121
$sub = $d.Range($sStart, $sEnd)
$sub.Font.Italic = $true

# paragraph 12
$sStart = $pStart11 + 22
$sEnd = $pStart11 + 57
$sub = $d.Range($sStart, $sEnd)
$sub.Font.Italic = $true

# paragraph 13
$sStart = $pStart12 + 25
$sEnd = $pStart12 + 61
$sub = $d.Range($sStart, $sEnd)
$sub.Font.Italic = $true

# ---- Update the header date from 2023-09-13 to 2023-09-15.
$d.Content.Find.Execute("2023-09-13", $true, $false, $false, $false, $false,
                        $true, 1, $false, "2023-09-15", 2) | Out-Null
